$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("evap")

# Remove the "treatment" column (column B) entirely; it only ever held the
# constant value "DMSO" for every row, so it carries no information.
$ws.Columns.Item(2).Delete()

# Leave the selection on the (now shifted) "time" column, matching the
# cursor position the author ended up with after deleting the column.
$ws.Columns.Item(2).Select()
